$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.01
$ws.Range("B7").Value = 5.962000000000001
$ws.Range("A8").Value = -22.366
$ws.Range("A10").Value = -21.598
$ws.Range("A12").Value = -21.551
$ws.Range("B15").Value = 5.313000000000001
$ws.Range("A18").Value = -21.974
$ws.Range("B18").Value = 5.683
$ws.Range("D18").Value = -8.608000000000001
$ws.Range("D19").Value = -7.912000000000001
$ws.Range("B20").Value = 7.489
$ws.Range("D27").Value = -8.182
$ws.Range("B29").Value = 5.023
$ws.Range("B30").Value = 6.001
$ws.Range("B31").Value = 6.401000000000001
$ws.Range("D31").Value = -8.353
$ws.Range("A37").Value = -20.272
$ws.Range("D38").Value = -8.567000000000002
$ws.Range("B40").Value = 8.962
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("D44").Value = -7.936999999999999
$ws.Range("D47").Value = -7.65
$ws.Range("B50").Value = 4.807
$ws.Range("A55").Value = -21.811
$ws.Range("D58").Value = -8.370999999999999
$ws.Range("D65").Value = -7.676
$ws.Range("A68").Value = -21.435
$ws.Range("B68").Value = 5.298
$ws.Range("D73").Value = -8.255000000000001
$ws.Range("B76").Value = 6.343000000000001
$ws.Range("A77").Value = -20.637
$ws.Range("A78").Value = -20.292
$ws.Range("A81").Value = -21.76
$ws.Range("A82").Value = -22.077
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858000000000001
$ws.Range("D90").Value = -8.237
$ws.Range("D94").Value = -7.302
$ws.Range("D95").Value = -7.718999999999999
$ws.Range("B96").Value = 6.468999999999999
$ws.Range("B98").Value = 5.828
$ws.Range("B101").Value = 8.194000000000001
$ws.Range("D101").Value = -8.16
$ws.Range("B102").Value = 7.306999999999999
